# object_rocket.xlsx survey-data edit
# - insert a "FAKE" column after Dislike (new col D)
# - insert a "Question_Pers_01" column after Next_Movie_To_Watch (new col K, before old Question_Rec_01)
# - rename old "Timestamp_page_1" column in place to "Question_Rec_20" (no longer a timestamp)
# - insert a "Timestamp_start_session" column right before Watchlist
# - replace the single data row (row 2) with the new survey response

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- structural column changes (left to right, so each Insert's target
#      letter is still valid at the moment it runs) ----

# 1) New column D: "FAKE"
$ws.Range("D1").EntireColumn.Insert()

# 2) New column K: "Question_Pers_01" (K currently holds old column J's shifted
#    content = Question_Rec_01, so inserting here pushes it to L as expected)
$ws.Range("K1").EntireColumn.Insert()

# 3) New column before Watchlist: "Timestamp_start_session"
#    (Watchlist is currently at AJ after the two inserts above)
$ws.Range("AJ1").EntireColumn.Insert()

# ---- header row (row 1) ----
$ws.Range("D1").Value = "FAKE"
$ws.Range("K1").Value = "Question_Pers_01"
# AE1 used to be "Timestamp_page_1"; it's now a plain recommendation question
$ws.Range("AE1").Value = "Question_Rec_20"
$ws.Range("AJ1").Value = "Timestamp_start_session"

# ---- data row (row 2) : new survey response ----
$ws.Range("A2").Value = "5e7bec2d9984da06eb63fb46"
$ws.Range("B2").Value = "26-40"
$ws.Range("C2").Value = "[22, 40]"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "test 6"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "male"
$ws.Range("H2").Value = "[9, 20, 36]"
$ws.Range("I2").Value = "[3, 29, 44]"
$ws.Range("J2").Value = "The Wolf of Wall Street "

# K2..AE2 hold numeric-looking *text* (answers "1".."7"), so force Text format
# first (otherwise Excel auto-coerces them to numbers), then strip the format
# back off so no stray cell style is left behind.
$textCells = @{
    "K2"  = "1"
    "L2"  = "7"
    "M2"  = "6"
    "N2"  = "5"
    "O2"  = "4"
    "P2"  = "3"
    "Q2"  = "2"
    "R2"  = "1"
    "S2"  = "1"
    "T2"  = "1"
    "U2"  = "1"
    "V2"  = "1"
    "W2"  = "1"
    "X2"  = "1"
    "Y2"  = "1"
    "Z2"  = "1"
    "AA2" = "1"
    "AB2" = "1"
    "AC2" = "1"
    "AD2" = "1"
    "AE2" = "2"
}
foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.ClearFormats()
}

$ws.Range("AF2").Value = 43915.9866840625
$ws.Range("AG2").Value = 43915.98673648148
$ws.Range("AH2").Value = 43915.98717037037
$ws.Range("AI2").Value = 43915.98718770834
$ws.Range("AJ2").Value = 43915.98652409722
$ws.Range("AK2").Value = "['The Prestige', 'The Wolf of Wall Street ']"
$ws.Range("AL2").Value = "5e7bec2d9984da06eb63fb46"
